{"js": "// Update two-digit multiplication equations to new values.\nconst replacements = [\n  [\"71\u00d725=1775\", \"43\u00d741=1763\"],\n  [\"35\u00d723=805\", \"40\u00d745=1800\"],\n  [\"72\u00d736=2592\", \"82\u00d791=7462\"],\n  [\"43\u00d722=946\", \"74\u00d773=5402\"],\n  [\"20\u00d782=1640\", \"29\u00d759=1711\"],\n  [\"78\u00d714=1092\", \"14\u00d761=854\"],\n  [\"54\u00d712=648\", \"53\u00d798=5194\"],\n  [\"84\u00d797=8148\", \"31\u00d728=868\"],\n  [\"46\u00d769=3174\", \"95\u00d798=9310\"],\n  [\"96\u00d731=2976\", \"72\u00d768=4896\"],\n  [\"24\u00d720=480\", \"84\u00d794=7896\"],\n  [\"60\u00d722=1320\", \"97\u00d764=6208\"],\n  [\"22\u00d722=484\", \"62\u00d771=4402\"],\n  [\"62\u00d729=1798\", \"30\u00d778=2340\"],\n  [\"83\u00d757=4731\", \"11\u00d760=660\"],\n  [\"66\u00d787=5742\", \"71\u00d753=3763\"],\n  [\"43\u00d785=3655\", \"89\u00d737=3293\"],\n  [\"84\u00d724=2016\", \"50\u00d771=3550\"],\n  [\"57\u00d779=4503\", \"79\u00d773=5767\"],\n  [\"31\u00d751=1581\", \"85\u00d717=1445\"],\n  [\"18\u00d789=1602\", \"65\u00d788=5720\"],\n  [\"49\u00d742=2058\", \"12\u00d794=1128\"],\n  [\"74\u00d714=1036\", \"80\u00d723=1840\"],\n  [\"50\u00d770=3500\", \"68\u00d748=3264\"],\n  [\"30\u00d774=2220\", \"60\u00d791=5460\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update two-digit multiplication equations to new values.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"71\u00d725=1775\", \"43\u00d741=1763\"),\n    @(\"35\u00d723=805\", \"40\u00d745=1800\"),\n    @(\"72\u00d736=2592\", \"82\u00d791=7462\"),\n    @(\"43\u00d722=946\", \"74\u00d773=5402\"),\n    @(\"20\u00d782=1640\", \"29\u00d759=1711\"),\n    @(\"78\u00d714=1092\", \"14\u00d761=854\"),\n    @(\"54\u00d712=648\", \"53\u00d798=5194\"),\n    @(\"84\u00d797=8148\", \"31\u00d728=868\"),\n    @(\"46\u00d769=3174\", \"95\u00d798=9310\"),\n    @(\"96\u00d731=2976\", \"72\u00d768=4896\"),\n    @(\"24\u00d720=480\", \"84\u00d794=7896\"),\n    @(\"60\u00d722=1320\", \"97\u00d764=6208\"),\n    @(\"22\u00d722=484\", \"62\u00d771=4402\"),\n    @(\"62\u00d729=1798\", \"30\u00d778=2340\"),\n    @(\"83\u00d757=4731\", \"11\u00d760=660\"),\n    @(\"66\u00d787=5742\", \"71\u00d753=3763\"),\n    @(\"43\u00d785=3655\", \"89\u00d737=3293\"),\n    @(\"84\u00d724=2016\", \"50\u00d771=3550\"),\n    @(\"57\u00d779=4503\", \"79\u00d773=5767\"),\n    @(\"31\u00d751=1581\", \"85\u00d717=1445\"),\n    @(\"18\u00d789=1602\", \"65\u00d788=5720\"),\n    @(\"49\u00d742=2058\", \"12\u00d794=1128\"),\n    @(\"74\u00d714=1036\", \"80\u00d723=1840\"),\n    @(\"50\u00d770=3500\", \"68\u00d748=3264\"),\n    @(\"30\u00d774=2220\", \"60\u00d791=5460\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: not found: $old\"\n    }\n}\n"}
